$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: Insert two new rows at position 3 (shifts old rows 3,4 down to 5,6)
$ws.Rows("3:4").Insert()

# Step 2: Update row 2 existing values per diff
$ws.Range("G2").Value = 2.45
$ws.Range("I2").Value = 3.1
$ws.Range("J2").Value = 3.25
$ws.Range("M2").Value = 1.1
$ws.Range("N2").Value = 7
$ws.Range("U2").Value = 2.1
$ws.Range("V2").Value = 1.67
$ws.Range("X2").Value = 10
$ws.Range("Z2").Value = 23
$ws.Range("AE2").Value = 19
$ws.Range("AF2").Value = 67
$ws.Range("AG2").Value = 7.5
$ws.Range("AH2").Value = 13
$ws.Range("AK2").Value = 29
$ws.Range("AU2").Value = 9
$ws.Range("AX2").Value = 19
$ws.Range("BA2").Value = 101

# Step 3: Fill new row 3 (Defensa y Justicia vs Dep. Riestra)
$ws.Range("A3").Value = "zyYw8Qbe"
$ws.Range("B3").Value = "18/11/2024"
$ws.Range("C3").Value = "19:15"
$ws.Range("D3").Value = "ARGENTINA - TORNEO BETANO"
$ws.Range("E3").Value = "Defensa y Justicia"
$ws.Range("F3").Value = "Dep. Riestra"
$ws.Range("G3").Value = 1.73
$ws.Range("H3").Value = 3.3
$ws.Range("I3").Value = 5.75
$ws.Range("J3").Value = 2.4
$ws.Range("K3").Value = 2.05
$ws.Range("L3").Value = 5.5
$ws.Range("M3").Value = 1.08
$ws.Range("N3").Value = 8
$ws.Range("O3").Value = 1.4
$ws.Range("P3").Value = 2.75
$ws.Range("Q3").Value = 2.3
$ws.Range("R3").Value = 1.6
$ws.Range("S3").Value = 1.5
$ws.Range("T3").Value = 2.5
$ws.Range("U3").Value = 2.1
$ws.Range("V3").Value = 1.67
$ws.Range("W3").Value = 5.5
$ws.Range("X3").Value = 7
$ws.Range("Y3").Value = 9
$ws.Range("Z3").Value = 13
$ws.Range("AA3").Value = 17
$ws.Range("AB3").Value = 34
$ws.Range("AC3").Value = 7
$ws.Range("AD3").Value = 6.5
$ws.Range("AE3").Value = 19
$ws.Range("AF3").Value = 67
$ws.Range("AG3").Value = 12
$ws.Range("AH3").Value = 26
$ws.Range("AI3").Value = 19
$ws.Range("AJ3").Value = 51
$ws.Range("AK3").Value = 51
$ws.Range("AL3").Value = 51
$ws.Range("AM3").Value = 501
$ws.Range("AN3").Value = 3.5
$ws.Range("AO3").Value = 9.5
$ws.Range("AP3").Value = 23
$ws.Range("AQ3").Value = 34
$ws.Range("AR3").Value = 67
$ws.Range("AS3").Value = 201
$ws.Range("AT3").Value = 2.5
$ws.Range("AU3").Value = 9.5
$ws.Range("AV3").Value = 67
$ws.Range("AW3").Value = 6.5
$ws.Range("AX3").Value = 29
$ws.Range("AY3").Value = 41
$ws.Range("AZ3").Value = 126
$ws.Range("BA3").Value = 151
$ws.Range("BB3").Value = 401
$ws.Range("BC3").Value = 126
$ws.Range("BD3").Value = 126

# Step 4: Fill new row 4 (Platense vs Godoy Cruz)
$ws.Range("A4").Value = "2NZV9nUr"
$ws.Range("B4").Value = "18/11/2024"
$ws.Range("C4").Value = "19:15"
$ws.Range("D4").Value = "ARGENTINA - TORNEO BETANO"
$ws.Range("E4").Value = "Platense"
$ws.Range("F4").Value = "Godoy Cruz"
$ws.Range("G4").Value = 2
$ws.Range("H4").Value = 3
$ws.Range("I4").Value = 4.5
$ws.Range("J4").Value = 2.88
$ws.Range("K4").Value = 1.8
$ws.Range("L4").Value = 5.5
$ws.Range("M4").Value = 1.14
$ws.Range("N4").Value = 5.5
$ws.Range("O4").Value = 1.73
$ws.Range("P4").Value = 2
$ws.Range("Q4").Value = 3.4
$ws.Range("R4").Value = 1.33
$ws.Range("S4").Value = 1.75
$ws.Range("T4").Value = 2.05
$ws.Range("U4").Value = 2.63
$ws.Range("V4").Value = 1.44
$ws.Range("W4").Value = 4.5
$ws.Range("X4").Value = 7.5
$ws.Range("Y4").Value = 11
$ws.Range("Z4").Value = 17
$ws.Range("AA4").Value = 23
$ws.Range("AB4").Value = 51
$ws.Range("AC4").Value = 5
$ws.Range("AD4").Value = 6.5
$ws.Range("AE4").Value = 26
$ws.Range("AF4").Value = 126
$ws.Range("AG4").Value = 7.5
$ws.Range("AH4").Value = 21
$ws.Range("AI4").Value = 19
$ws.Range("AJ4").Value = 51
$ws.Range("AK4").Value = 51
$ws.Range("AL4").Value = 67
$ws.Range("AM4").Value = 201
$ws.Range("AN4").Value = 3.6
$ws.Range("AO4").Value = 13
$ws.Range("AP4").Value = 34
$ws.Range("AQ4").Value = 51
$ws.Range("AR4").Value = 101
$ws.Range("AS4").Value = 451
$ws.Range("AT4").Value = 2
$ws.Range("AU4").Value = 11
$ws.Range("AV4").Value = 101
$ws.Range("AW4").Value = 6
$ws.Range("AX4").Value = 34
$ws.Range("AY4").Value = 51
$ws.Range("AZ4").Value = 126
$ws.Range("BA4").Value = 251
$ws.Range("BB4").Value = 501
$ws.Range("BC4").Value = 126
$ws.Range("BD4").Value = 126
